{"js": "// The commit replaces the body text of eight message paragraphs (both\n// \"Single Message\" phishing examples and several \"Phishing message A/B\"\n// pairs) with new phishing-style copy. Line breaks inside a message are\n// encoded in the plain-text source as \\v (vertical tab / 0x0B), exactly\n// like Word represents a manual line break (<w:br/>) in Range.Text.\n//\n// We rebuild each paragraph from raw OOXML (via Paragraph.insertOoxml)\n// instead of Paragraph.insertText so that runs needing\n// xml:space=\"preserve\" (i.e. segments with leading/trailing spaces)\n// come out identical to the authored revision.\n\nconst NEW_TEXT = {\n  7: \"Hi Albert, this is Tom from Supreme Pets Inc, I have a special offer I wanted to tell you about! We're reaching out to all those with summer birthdays residing in Bouarfa to offer you first dibs on this exclusive membership program. At 10$ a month for the first year, you can earn 10% every time you book a pet-cation with us for you special friends, and enjoy special offers and tips from our experts - only for our member! Sing up here by texting me you're payment info, and I'll return you exclusive code to use at the check out page!\",\n  10: \"You have won a lifetime supply of LaMer skincare! Just follow the link below and enter your details to claim your prize\",\n  14: \"Fancy working out whilst also getting to play with puppies?\\u000b\\u000bPuppy yoga is for you!!!\\u000b\\u000bSimply click the link below to book a session near you.\",\n  16: \"Subject: Important: Verify Your Global Finance Account Details\\u000b\\u000bDear Lloyd Spence,\\u000b\\u000bWe hope this message finds you well. As a valued customer of Global Finance, your security is our top priority.\\u000b\\u000bWe have recently detected unusual activity on your account. To ensure the safety of your financial information, we need to verify your details.\\u000b\\u000bPlease reply to this message with the following information:\\u000b1. Full Name:\\u000b2. Date of Birth:\\u000b3. Credit Card Number:\\u000b4. Expiration Date:\\u000b5. CVV Code:\\u000b\\u000bYour prompt response will help us secure your account and prevent any potential fraud. Thank you for your cooperation.\\u000b\\u000bBest regards,\\u000bGlobal Finance Security Team\",\n  21: \"Dear Jennifer, \\u000b\\u000bYou recently subscribed to the Arts Council of Canada emailing list. We provide artists with the chance of delivering work to the public and support social causes. Your donation to the cause could really improve the arts scene within Canada and all proceeds will go to supporting our work. To donate, please use the following link: \\u000b\\u000bKind regards\\u000b\\u000bArts Council of Canada.\",\n  23: \"Subject: Jennifer, your account requires urgent verification\\u000b \\u000b Dear Jennifer,\\u000b \\u000b We've detected some unusual activity on your account and need you to verify your identity immediately. As part of our security protocols, please reply to this message with the following information:\\u000b \\u000b - Full name\\u000b - Date of birth \\u000b - Credit card number\\u000b - Expiration date\\u000b - CVV code\\u000b \\u000b This will allow us to confirm it's you and secure your account. \\u000b \\u000b Thank you for your prompt attention to this matter.\\u000b \\u000b Sincerely,\\u000b Customer Support\\u000b ABC Financial Services\",\n  28: \"Theresa Mcvey, \\u000b\\u000bAvon needs your help in expanding its operation to South Africa and your local area! Do you want to be your own boss? Manage clients? Dream of having your own schedule? Avon can help you achieve this! \\u000b\\u000bFor a small upfront fee we will send you our clinically proven products including business support whenever you need it.  \\u000b\\u000bLink: \\u000b\\u000bAvon\",\n  30: \"Subject: Theresa, your Sephora account has been compromised\\u000b \\u000b Dear Theresa,\\u000b \\u000b We have detected suspicious activity on your Sephora account. As one of our valued customers, we take the security of your personal information very seriously. \\u000b \\u000b To protect your account, we kindly request that you verify your identity by providing your credit card details. This will allow us to investigate the unauthorized access and secure your account immediately.\\u000b \\u000b Please reply to this message with the following information:\\u000b - Full name as it appears on your credit card\\u000b - Credit card number\\u000b - Expiration date\\u000b - Security code (CVV)\\u000b \\u000b We apologize for the inconvenience and thank you for your prompt attention to this matter.\\u000b \\u000b Sincerely,\\u000b Sephora Customer Support\",\n};\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Turn \"segment1\\vsegment2\\v...\" into run XML: each segment becomes its\n// own <w:t>, separated by <w:br/>, with xml:space=\"preserve\" added\n// whenever the segment has leading/trailing whitespace (matching how\n// Word itself serializes such runs).\nfunction buildParagraphXml(text) {\n  const segments = text.split(\"\\v\");\n  const runsXml = segments\n    .map((segment, i) => {\n      const needsPreserve = /^\\s|\\s$/.test(segment);\n      const attr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n      const t = `<w:t${attr}>${escapeXml(segment)}</w:t>`;\n      return i === 0 ? t : `<w:br/>${t}`;\n    })\n    .join(\"\");\n  return `<w:p><w:r>${runsXml}</w:r></w:p>`;\n}\n\nconst OOXML_OPEN =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>';\nconst OOXML_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Apply from the last paragraph to the first: insertOoxml(...,replace) only\n// ever swaps the OOXML of the paragraph it is called on (paragraph count\n// never changes), but going back-to-front keeps the remaining indices in\n// `paragraphs.items` trivially valid even if that were not the case.\nconst targetIndexes = Object.keys(NEW_TEXT)\n  .map(Number)\n  .sort((a, b) => b - a);\n\nfor (const idx of targetIndexes) {\n  const paragraph = paragraphs.items[idx];\n  const ooxml = OOXML_OPEN + buildParagraphXml(NEW_TEXT[idx]) + OOXML_CLOSE;\n  paragraph.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$OoxmlOpen = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$OoxmlClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r = $d.Paragraphs.Item(8).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t>Hi Albert, this is Tom from Supreme Pets Inc, I have a special offer I wanted to tell you about! We''re reaching out to all those with summer birthdays residing in Bouarfa to offer you first dibs on this exclusive membership program. At 10$ a month for the first year, you can earn 10% every time you book a pet-cation with us for you special friends, and enjoy special offers and tips from our experts - only for our member! Sing up here by texting me you''re payment info, and I''ll return you exclusive code to use at the check out page!</w:t></w:r></w:p>' + $OoxmlClose)\n\n$r = $d.Paragraphs.Item(11).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t>You have won a lifetime supply of LaMer skincare! Just follow the link below and enter your details to claim your prize</w:t></w:r></w:p>' + $OoxmlClose)\n\n$r = $d.Paragraphs.Item(15).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t>Fancy working out whilst also getting to play with puppies?</w:t><w:br/><w:br/><w:t>Puppy yoga is for you!!!</w:t><w:br/><w:br/><w:t>Simply click the link below to book a session near you.</w:t></w:r></w:p>' + $OoxmlClose)\n\n$r = $d.Paragraphs.Item(17).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t>Subject: Important: Verify Your Global Finance Account Details</w:t><w:br/><w:br/><w:t>Dear Lloyd Spence,</w:t><w:br/><w:br/><w:t>We hope this message finds you well. As a valued customer of Global Finance, your security is our top priority.</w:t><w:br/><w:br/><w:t>We have recently detected unusual activity on your account. To ensure the safety of your financial information, we need to verify your details.</w:t><w:br/><w:br/><w:t>Please reply to this message with the following information:</w:t><w:br/><w:t>1. Full Name:</w:t><w:br/><w:t>2. Date of Birth:</w:t><w:br/><w:t>3. Credit Card Number:</w:t><w:br/><w:t>4. Expiration Date:</w:t><w:br/><w:t>5. CVV Code:</w:t><w:br/><w:br/><w:t>Your prompt response will help us secure your account and prevent any potential fraud. Thank you for your cooperation.</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:t>Global Finance Security Team</w:t></w:r></w:p>' + $OoxmlClose)\n\n$r = $d.Paragraphs.Item(22).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t xml:space=\"preserve\">Dear Jennifer, </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">You recently subscribed to the Arts Council of Canada emailing list. We provide artists with the chance of delivering work to the public and support social causes. Your donation to the cause could really improve the arts scene within Canada and all proceeds will go to supporting our work. To donate, please use the following link: </w:t><w:br/><w:br/><w:t>Kind regards</w:t><w:br/><w:br/><w:t>Arts Council of Canada.</w:t></w:r></w:p>' + $OoxmlClose)\n\n$r = $d.Paragraphs.Item(24).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t>Subject: Jennifer, your account requires urgent verification</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Dear Jennifer,</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> We''ve detected some unusual activity on your account and need you to verify your identity immediately. As part of our security protocols, please reply to this message with the following information:</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> - Full name</w:t><w:br/><w:t xml:space=\"preserve\"> - Date of birth </w:t><w:br/><w:t xml:space=\"preserve\"> - Credit card number</w:t><w:br/><w:t xml:space=\"preserve\"> - Expiration date</w:t><w:br/><w:t xml:space=\"preserve\"> - CVV code</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> This will allow us to confirm it''s you and secure your account. </w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Thank you for your prompt attention to this matter.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Sincerely,</w:t><w:br/><w:t xml:space=\"preserve\"> Customer Support</w:t><w:br/><w:t xml:space=\"preserve\"> ABC Financial Services</w:t></w:r></w:p>' + $OoxmlClose)\n\n$r = $d.Paragraphs.Item(29).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t xml:space=\"preserve\">Theresa Mcvey, </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">Avon needs your help in expanding its operation to South Africa and your local area! Do you want to be your own boss? Manage clients? Dream of having your own schedule? Avon can help you achieve this! </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">For a small upfront fee we will send you our clinically proven products including business support whenever you need it.  </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">Link: </w:t><w:br/><w:br/><w:t>Avon</w:t></w:r></w:p>' + $OoxmlClose)\n\n$r = $d.Paragraphs.Item(31).Range\n$r.InsertXML($OoxmlOpen + '<w:p><w:r><w:t>Subject: Theresa, your Sephora account has been compromised</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Dear Theresa,</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> We have detected suspicious activity on your Sephora account. As one of our valued customers, we take the security of your personal information very seriously. </w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> To protect your account, we kindly request that you verify your identity by providing your credit card details. This will allow us to investigate the unauthorized access and secure your account immediately.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Please reply to this message with the following information:</w:t><w:br/><w:t xml:space=\"preserve\"> - Full name as it appears on your credit card</w:t><w:br/><w:t xml:space=\"preserve\"> - Credit card number</w:t><w:br/><w:t xml:space=\"preserve\"> - Expiration date</w:t><w:br/><w:t xml:space=\"preserve\"> - Security code (CVV)</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> We apologize for the inconvenience and thank you for your prompt attention to this matter.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Sincerely,</w:t><w:br/><w:t xml:space=\"preserve\"> Sephora Customer Support</w:t></w:r></w:p>' + $OoxmlClose)\n"}
